# Mise à jour de l'application
# Ajout de 3 nouvelles dates d'entrainement (colonnes CH, CI, CJ) avec les
# présences/absences correspondantes pour chaque joueur.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Nouvelles dates en ligne 1 -----------------------------------------
$ws.Cells.Item(1, 86).Value = 45988   # CH1
$ws.Cells.Item(1, 87).Value = 45989   # CI1
$ws.Cells.Item(1, 88).Value = 45993   # CJ1

# --- 2. Valeurs de présence pour chaque joueur (lignes 2 à 29) -------------
# Lignes avec la même valeur dans CH, CI et CJ ("P" = présent)
$ws.Range("CH2:CJ4").Value   = "P"
$ws.Range("CH7:CJ11").Value  = "P"
$ws.Range("CH14:CJ20").Value = "P"
$ws.Range("CH22:CJ22").Value = "P"
$ws.Range("CH24:CJ24").Value = "P"
$ws.Range("CH27:CJ29").Value = "P"

# Lignes avec "RH" (réserve/hors groupe) sur les 3 colonnes
$ws.Range("CH6:CJ6").Value = "RH"

# Lignes avec "B" (blessure) sur les 3 colonnes
$ws.Range("CH13:CJ13").Value = "B"
$ws.Range("CH25:CJ25").Value = "B"

# Ligne 5 : RH, RH, P
$ws.Cells.Item(5, 86).Value = "RH"
$ws.Cells.Item(5, 87).Value = "RH"
$ws.Cells.Item(5, 88).Value = "P"

# Ligne 23 : P, R, RH
$ws.Cells.Item(23, 86).Value = "P"
$ws.Cells.Item(23, 87).Value = "R"
$ws.Cells.Item(23, 88).Value = "RH"

# Ligne 26 : P, B, P
$ws.Cells.Item(26, 86).Value = "P"
$ws.Cells.Item(26, 87).Value = "B"
$ws.Cells.Item(26, 88).Value = "P"

# Ligne 21 : le joueur n'est plus suivi depuis plusieurs colonnes déjà :
# les cellules restent vides (aucune valeur à écrire).
# Ligne 12 : le joueur n'a plus de colonnes au-delà de AX, on ne touche pas
# aux colonnes CH/CI/CJ pour cette ligne.

# --- 3. Recopie de la mise en forme (style) depuis la colonne CG ----------
# (copie des formats uniquement, sans écraser les valeurs déjà saisies
#  ci-dessus ; la ligne 12 est volontairement exclue car elle ne va pas
#  jusqu'à la colonne CG)
$ws.Range("CG1").Copy()
$ws.Range("CH1:CJ1").PasteSpecial(-4122)

$ws.Range("CG2:CG11").Copy()
$ws.Range("CH2:CH11").PasteSpecial(-4122)
$ws.Range("CI2:CI11").PasteSpecial(-4122)
$ws.Range("CJ2:CJ11").PasteSpecial(-4122)

$ws.Range("CG13:CG29").Copy()
$ws.Range("CH13:CH29").PasteSpecial(-4122)
$ws.Range("CI13:CI29").PasteSpecial(-4122)
$ws.Range("CJ13:CJ29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 4. Sélection active telle qu'elle était à l'enregistrement -----------
$ws.Range("CL24").Select()
